$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.84
$ws.Range("H2").Value = 2.64
$ws.Range("I2").Value = 2.66
$ws.Range("J2").Value = 3.65
$ws.Range("K2").Value = 3.7
$ws.Range("P2").Value = 2.16
$ws.Range("V2").Value = 1.6
$ws.Range("Z2").Value = 19.5
$ws.Range("AM2").Value = 75
$ws.Range("AO2").Value = 19.5
$ws.Range("F3").Value = 16.5
$ws.Range("H3").Value = 1.15
$ws.Range("K3").Value = 16
$ws.Range("V3").Value = 5.8
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AC3").Value = 980
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AH3").Value = 980
$ws.Range("AI3").Value = 980
$ws.Range("F4").Value = 2.64
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 2.46
$ws.Range("I4").Value = 2.74
$ws.Range("J4").Value = 3.55
$ws.Range("O4").Value = 1.24
$ws.Range("Q4").Value = 1.72
$ws.Range("S4").Value = 2.84
$ws.Range("T4").Value = 1.61
$ws.Range("V4").Value = 1.57
$ws.Range("W4").Value = 1.5
$ws.Range("X4").Value = 980
$ws.Range("Y4").Value = 15.5
$ws.Range("AA4").Value = 980
$ws.Range("AB4").Value = 16.5
$ws.Range("AC4").Value = 10.5
$ws.Range("AE4").Value = 980
$ws.Range("AF4").Value = 980
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 19.5
$ws.Range("AI4").Value = 980
$ws.Range("AJ4").Value = 980
$ws.Range("AK4").Value = 980
$ws.Range("AL4").Value = 980
$ws.Range("AM4").Value = 1000
$ws.Range("AN4").Value = 980
$ws.Range("AO4").Value = 23
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 4.8
$ws.Range("H5").Value = 1.78
$ws.Range("I5").Value = 1.88
$ws.Range("M5").Value = 1.02
$ws.Range("N5").Value = 4.9
$ws.Range("O5").Value = 1.2
$ws.Range("Q5").Value = 1.6
$ws.Range("S5").Value = 2.48
$ws.Range("T5").Value = 1.05
$ws.Range("U5").Value = 1.05
$ws.Range("V5").Value = 2.04
$ws.Range("W5").Value = 1.26
$ws.Range("X5").Value = 28
$ws.Range("Y5").Value = 15
$ws.Range("Z5").Value = 16.5
$ws.Range("AA5").Value = 25
$ws.Range("AG5").Value = 22
$ws.Range("AJ5").Value = 100
$ws.Range("AM5").Value = 85
$ws.Range("AN5").Value = 44
$ws.Range("AO5").Value = 10.5
$ws.Range("H6").Value = 2.18
$ws.Range("I6").Value = 2.2
$ws.Range("L6").Value = 1.39
$ws.Range("N6").Value = 4.1
$ws.Range("O6").Value = 1.3
$ws.Range("P6").Value = 2.06
$ws.Range("Q6").Value = 1.91
$ws.Range("R6").Value = 1.4
$ws.Range("S6").Value = 3.35
$ws.Range("T6").Value = 1.74
$ws.Range("U6").Value = 2.24
$ws.Range("V6").Value = 1.83
$ws.Range("X6").Value = 14.5
$ws.Range("Z6").Value = 13.5
$ws.Range("AB6").Value = 15
$ws.Range("AI6").Value = 34
$ws.Range("AM6").Value = 85
$ws.Range("AO6").Value = 15
$ws.Range("F7").Value = 3.05
$ws.Range("J7").Value = 2.76
$ws.Range("N7").Value = 2.92
$ws.Range("P7").Value = 1.65
$ws.Range("T7").Value = 1.89
$ws.Range("Z7").Value = 17
$ws.Range("AA7").Value = 980
$ws.Range("AE7").Value = 980
$ws.Range("AF7").Value = 980
$ws.Range("AG7").Value = 19
$ws.Range("AH7").Value = 980
$ws.Range("AJ7").Value = 1000
$ws.Range("AL7").Value = 1000
$ws.Range("AM7").Value = 1000
$ws.Range("AO7").Value = 980
$ws.Range("G8").Value = 2.52
$ws.Range("H8").Value = 3
$ws.Range("L8").Value = 1.42
$ws.Range("W8").Value = 1.65
$ws.Range("G9").Value = 5.4
$ws.Range("I9").Value = 3.5
$ws.Range("J9").Value = 1.01
$ws.Range("V9").Value = 1.41
$ws.Range("W9").Value = 1.29
$ws.Range("F11").Value = 2.94
$ws.Range("G11").Value = 3.8
$ws.Range("I11").Value = 2.88
$ws.Range("J11").Value = 2.88
$ws.Range("N11").Value = 2.86
$ws.Range("O11").Value = 1.45
$ws.Range("Q11").Value = 2.32
$ws.Range("S11").Value = 3.95
$ws.Range("T11").Value = 1.94
$ws.Range("V11").Value = 1.58
$ws.Range("W11").Value = 1.4
$ws.Range("Y11").Value = 1000
$ws.Range("AB11").Value = 1000
$ws.Range("AD11").Value = 1000
$ws.Range("AE11").Value = 980
$ws.Range("AG11").Value = 1000
$ws.Range("F12").Value = 1.58
$ws.Range("G12").Value = 1.66
$ws.Range("H12").Value = 5.2
$ws.Range("I12").Value = 6.4
$ws.Range("J12").Value = 4.3
$ws.Range("K12").Value = 5
$ws.Range("P12").Value = 2.32
$ws.Range("R12").Value = 1.52
$ws.Range("S12").Value = 2.54
$ws.Range("T12").Value = 1.74
$ws.Range("U12").Value = 2.16
$ws.Range("V12").Value = 1.19
$ws.Range("W12").Value = 2.5
$ws.Range("AF12").Value = 970
$ws.Range("AL12").Value = 970
$ws.Range("H13").Value = 2.4
$ws.Range("P13").Value = 1.96
$ws.Range("Y13").Value = 1000
$ws.Range("AA13").Value = 980
$ws.Range("AF13").Value = 980
$ws.Range("AH13").Value = 1000
$ws.Range("AO13").Value = 1000
$ws.Range("R14").Value = 1.41
$ws.Range("T14").Value = 1.78
$ws.Range("AL14").Value = 40
$ws.Range("U15").Value = 2.12
$ws.Range("Z15").Value = 970
$ws.Range("AE15").Value = 970
$ws.Range("AI15").Value = 970
$ws.Range("AJ15").Value = 970
$ws.Range("AK15").Value = 970
$ws.Range("AL15").Value = 970
$ws.Range("AN15").Value = 970
$ws.Range("AO15").Value = 970
$ws.Range("N16").Value = 3.6
$ws.Range("R16").Value = 1.35
$ws.Range("S16").Value = 3.25
$ws.Range("AB16").Value = 11.5
$ws.Range("AI16").Value = 50
$ws.Range("AM16").Value = 100
$ws.Range("F17").Value = 1.71
$ws.Range("G17").Value = 2.18
$ws.Range("H17").Value = 3.2
$ws.Range("J17").Value = 3.2
$ws.Range("K17").Value = 9.6
$ws.Range("P17").Value = 1.78
$ws.Range("Q17").Value = 1.01
$ws.Range("S17").Value = 1.05
$ws.Range("W17").Value = 1.84
$ws.Range("H18").Value = 3.85
$ws.Range("W18").Value = 1.91
$ws.Range("H19").Value = 1.61
$ws.Range("I19").Value = 1.63
$ws.Range("P19").Value = 2.68
$ws.Range("S19").Value = 2.42
$ws.Range("H20").Value = 1.97
$ws.Range("G21").Value = 2.02
$ws.Range("H21").Value = 4.8
$ws.Range("J21").Value = 3.3
$ws.Range("W21").Value = 1.98
$ws.Range("F22").Value = 3.55
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 2.16
$ws.Range("I22").Value = 2.46
$ws.Range("N22").Value = 2.42
$ws.Range("S22").Value = 4.4
$ws.Range("V22").Value = 1.69
$ws.Range("W22").Value = 1.36
$ws.Range("X22").Value = 12
$ws.Range("Y22").Value = 9
$ws.Range("AB22").Value = 13
$ws.Range("AD22").Value = 970
$ws.Range("AJ22").Value = 110
$ws.Range("AL22").Value = 95
$ws.Range("AM22").Value = 210
$ws.Range("F23").Value = 2.22
$ws.Range("G23").Value = 2.42
$ws.Range("H23").Value = 3.55
$ws.Range("J23").Value = 3.1
$ws.Range("K23").Value = 3.45
$ws.Range("N23").Value = 2.62
$ws.Range("O23").Value = 1.52
$ws.Range("P23").Value = 1.52
$ws.Range("Q23").Value = 2.3
$ws.Range("R23").Value = 1.19
$ws.Range("S23").Value = 2.68
$ws.Range("T23").Value = 2.12
$ws.Range("U23").Value = 1.75
$ws.Range("W23").Value = 1.7
$ws.Range("X23").Value = 11
$ws.Range("Y23").Value = 12.5
$ws.Range("AB23").Value = 8.800000000000001
$ws.Range("AC23").Value = 8.800000000000001
